# Update Leave Card 4/11/2023 10:05 PM
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$tbl = $ws.ListObjects.Item("Table1")

# Insert a new blank row right below the existing 3/1/23 (row 132) entry.
# This shifts every subsequent row (old 133..177) down by one (new 134..178),
# including the specially-styled trailing "blank" row of the table which
# lands on row 178 automatically.
$ws.Rows.Item(133).Insert()

# The new row was created with generic default formatting; pull the
# (correct) formatting that was just pushed down from the old row 133
# into row 134, back up onto the blank row 133.
$ws.Range("A134:K134").Copy() | Out-Null
$ws.Range("A133:K133").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Grow the table so it covers the newly inserted row at the bottom (178)
# as well as the freshly inserted row 133 (already inside A8:K177).
$tbl.Resize($ws.Range("A8:K178"))

# Restore / set the calculated "EARNED " helper-column formula for the two
# rows that lost it: the brand-new row 133 (PasteSpecial above only copied
# formats, not formulas) and the new last row 178 -- structured references
# don't auto-fill there because the row only became part of the table after
# the Resize above.
$ws.Range("G133").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'
$ws.Range("G178").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'

# --- Row 132: SL(6-0-0), 3/16-23/2023 -------------------------------------
$ws.Range("B132").Value = "SL(6-0-0)"
$ws.Range("C132").Value = 1.25
$ws.Range("H132").Value = 6
$ws.Range("K132").Value = "3/16-23/2023"

# --- Row 133 (new): SL(9-0-0), 3/24 - 4/5/2023 ----------------------------
$ws.Range("B133").Value = "SL(9-0-0)"
$ws.Range("H133").Value = 9
$ws.Range("K133").Value = "3/24 - 4/5/2023"

# Match the author's last recorded selection.
$ws.Range("I134").Select()
